# Add the new "Balance Track" entry row (row 3) to Sheet1:
#   A3 = "4007.37 Euro"  (new shared string, left-aligned, #,##0.00 number format)
#   B3 = 45463           (date 2024-06-20, same date format as B2)
# and move the active selection to K18, matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New balance text in column A, row 3
$ws.Range("A3").Value = "4007.37 Euro"

# New date in column B, row 3 - reuse the same date number format as B2
$ws.Range("B3").Value = 45463
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat()

# A3 gets its own style: numeric format "#,##0.00" with left horizontal alignment
$ws.Range("A3").NumberFormat = "#,##0.00"
$ws.Range("A3").HorizontalAlignment = -4131   # xlLeft

# Move/select K18 as the last active cell, like in the saved workbook
[void]$ws.Range("K18").Select()
